$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two rows that were removed entirely (RM 232 at row 26, SC 92 at row 28).
# Delete from the bottom up so row indices for the earlier deletion stay valid.
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(26).Delete()

# Final data (after the two row deletions, rows 2-33) as it should read after the edit.
$data = @(
    @("RM 2",    -19.7, 14.9,  -13.5, -7.2,  18.03),
    @("RM 8",    -19.7, 11.2,  -14.2, $null, $null),
    @("RM 9",    -18.7, $null, -15.4, -6.4,  17.97),
    @("RM 14",   -19.5, 12.3,  -14.4, $null, 17.66),
    @("RM 21",   -19.8, $null, -14.2, -5.7,  16.43),
    @("RM 32",   -19.5, 15.0,  -13.8, -7.1,  17.24),
    @("RM 38",   -19.9, 15.5,  -13.9, $null, $null),
    @("RM 42",   -20.6, 10.5,  -14.5, -6.8,  17.26),
    @("RM 52 a", -19.8, 11.5,  -14.7, -6.1,  16.43),
    @("RM 58",   -20.8, 11.4,  -15.5, -7.9,  17.65),
    @("RM 81",   $null, 12.5,  -14.1, -5.3,  17.45),
    @("RM 88",   -19.9, 12.5,  -13.9, -5.3,  17.1),
    @("RM 90",   -19.6, $null, -13.1, -5.4,  17.76),
    @("RM 95",   -19.1, 12.5,  -15.2, -8.4,  16.2),
    @("RM 103",  -19.5, 13.5,  -15.3, -6.9,  17.34),
    @("RM 116",  -19.4, 11.2,  -14.7, -7.3,  17.78),
    @("RM 120",  -19.6, 11.5,  -15.2, -8.5,  18.35),
    @("RM 125",  -20.6, 13.2,  -15.5, $null, $null),
    @("RM 134",  -19.0, 12.5,  -14.0, -7.2,  17.73),
    @("RM 135",  -18.9, 12.7,  -14.3, -8.7,  16.58),
    @("RM 138",  -19.3, $null, -15.4, $null, 16.81),
    @("RM 140",  -19.5, $null, -13.9, -7.0,  16.48),
    @("RM 142a", -17.7, $null, -13.9, -8.1,  16.78),
    @("RM 145",  -19.5, 10.7,  -15.5, -7.1,  16.6),
    @("SC 5",    -20.2, 10.8,  -13.8, -5.0,  17.38),
    @("SC 101",  $null, 10.0,  -14.6, $null, 17.0),
    @("SC 105",  -19.6, 11.1,  -13.7, -5.9,  17.44),
    @("SC 119",  -19.5, 11.2,  -13.0, $null, 18.06),
    @("SC 120",  -19.7, 11.4,  -13.6, -5.7,  16.89),
    @("SC 132",  -18.8, 15.3,  -13.7, -8.1,  $null),
    @("SC 193",  $null, 10.5,  -14.7, -6.4,  17.39),
    @("SC 232",  -19.5, 10.4,  -14.1, -10.7, 17.53)
)

$r = 2
foreach ($row in $data) {
    for ($c = 0; $c -lt 6; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r = $r + 1
}
